$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, J, K, L, M, P across rows 2..6.
# Each row takes on the values previously held by the row above it
# (cyclically: row 2 takes what row 6 had).
$newValues = @{
    2 = @{ D = 44608; J = 120; K = 600;  L = 650;  M = 625;  P = 625  }
    3 = @{ D = 44532; J = 60;  K = 2000; L = 2200; M = 2100; P = 2100 }
    4 = @{ D = 44610; J = 100; K = 600;  L = 650;  M = 625;  P = 625  }
    5 = @{ D = 44624; J = 120; K = 650;  L = 700;  M = 675;  P = 675  }
    6 = @{ D = 44533; J = 100; K = 2000; L = 2200; M = 2100; P = 2100 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
